$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers - reorder columns
$ws.Range("A1").Value = "living_rooms_1"
$ws.Range("B1").Value = "bedrooms_1"
$ws.Range("C1").Value = "kitchens_1"
$ws.Range("D1").Value = "living_rooms_2"
$ws.Range("E1").Value = "bedrooms_2"
$ws.Range("F1").Value = "kitchens_2"

# Row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0

# Row 3
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1

# Row 4
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0

# Row 5
$ws.Range("A5").Value = 0
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0

# Row 6
$ws.Range("A6").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0

# Row 7
$ws.Range("A7").Value = 0
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
